# Opportunity creation flow fix: write the org-linkage test value into the
# OrgData worksheet (2nd row, column A) and leave the selection where the
# recorded macro session left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OrgData")

$ws.Range("A2").Value = "hamas_258"

$ws.Activate()
$ws.Range("I9").Select()
